# Updates prices / 1h-volume % / coin name & link columns on the "cryptos" sheet
# to match the refreshed Coinranking snapshot (GitHub Actions data sync).
# Values are written with a leading apostrophe where needed so Excel stores them
# as text (matching the sheet's existing inline-string cells) instead of coercing
# numeric-looking strings (e.g. "300.05", "-6.74%") into Number/Percentage cells.
# Style is reset to "Normal" right after so no stray quote-prefix formatting sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''300.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''-6.74%'
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = '''35.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''-3.10%'
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = '''4.982'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''-2.85%'
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = '''0.07926'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''-1.76%'
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = '''1.910'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''-11.31%'
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = '''7.736'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''-4.17%'
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '''2.925'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''4.49%'
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '''0.9244'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-0.51%'
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '''0.1127'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''12.74%'
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '''0.1827'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''-3.03%'
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = '''0.09283'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''0.37%'
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = '''0.03530'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''-1.29%'
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '''0.09881'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''-0.56%'
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '''0.001397'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''-2.15%'
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''0.005742'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''1.57%'
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''3.495'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''1.10%'
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = '''4.014'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''-2.78%'
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("E19").Value = '''2.02%'
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("E20").Value = '''-1.58%'
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = '''5.037'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''-0.77%'
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = '''0.2399'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''8.96%'
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").Value = '''0.04501'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''-2.22%'
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = '''0.001214'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''-2.29%'
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = '''0.004579'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''-3.49%'
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").Value = '''0.0001250'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''-3.91%'
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("E27").Value = '''-6.83%'
$ws.Range("E27").Style = "Normal"
# Row 39
$ws.Range("D39").Value = '''0.01881'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''-4.37%'
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").Value = '''0.04686'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''-6.00%'
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = '''0.007602'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''-3.01%'
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").Value = '''0.009565'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''24.37%'
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("E43").Value = '''-5.68%'
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = '''0.002120'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''2.52%'
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = '''0.01111'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-5.95%'
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").Value = '''0.00006022'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''-5.60%'
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("E47").Value = '''-0.03%'
$ws.Range("E47").Style = "Normal"
# Row 49
$ws.Range("E49").Value = '''-31.36%'
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").Value = '''0.00002100'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''-0.03%'
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").Value = '''0.0002000'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''-0.03%'
$ws.Range("E51").Style = "Normal"
